# daily auto push: append today's row (2025-10-03, Fri) to the tracking sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 54

# Column A holds the date as literal text (e.g. "2025/09/22"), not a real
# date value. Force the cell to Text format before assigning the string so
# Excel doesn't auto-convert "2025/10/03" into a date serial number, then
# clear the formatting back to the sheet's default (no explicit style),
# matching the rest of the data rows.
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "2025/10/03"
$ws.Cells.Item($newRow, 1).ClearFormats()

$ws.Cells.Item($newRow, 2).Value = "金"
$ws.Cells.Item($newRow, 3).Value = 2
$ws.Cells.Item($newRow, 4).Value = 201
